# Corrected ICDC Breed 1-14 scripts:
#  - The "startup" sheet's CasesTab/SamplesTab/FilesTab query cells (B2:B4)
#    had their shared-string order corrected and minor whitespace issues
#    (stray spaces before commas) cleaned up in the Cases query and the
#    Samples query text.
#  - The active selection on the "startup" sheet moved from B4 to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Corrected "CasesTab" query (row 2, column B) - removed stray spaces before
# some commas in the RETURN clause.
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Akita']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID`,
        coalesce(s.clinical_study_designation, '') AS `Study Code`,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease`,
        coalesce(demo.patient_age_at_enrollment, '') AS Age,
        coalesce(demo.sex, '') AS Sex,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
        coalesce(co.cohort_description, '') AS `Cohort`

'@

# Corrected "SamplesTab" query (row 3, column B) - removed stray spaces
# before some commas in the RETURN clause.
$samplesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE demo.breed IN ['Akita'] 
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed,
        coalesce(diag.disease_term,'') AS Diagnosis, 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
'@

# "FilesTab" query (row 4, column B) - unchanged text content.
$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Akita']  
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B2").Value2 = $casesQuery
$ws.Range("B3").Value2 = $samplesQuery
$ws.Range("B4").Value2 = $filesQuery

# Move the active cell/selection from B4 to B2 and scroll the window so
# row 2 is the top-left visible row.
$ws.Activate()
[void]$ws.Range("B2").Select()
